$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# "Bento object repository revisited": the FilesTab Cypher query (cell B4)
# drops the `File Type` and `Breed` output columns from its RETURN clause.
$newQuery = @"
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
 MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
 MATCH (samp:sample)-->(c) 
 WHERE f.file_type IN ["RNA Sequence File"]  
WITH DISTINCT f, parent, c, demo, diag, s
RETURN coalesce(f.file_name, '') AS ``File Name``, 
        coalesce(labels(parent)[0], '') AS ``Association``,
        coalesce(f.file_description, '') AS ``Description``,
        coalesce(f.file_format, '') AS ``Format``,
        coalesce(f.file_size, '') AS ``Size``,
        coalesce(c.case_id, '') AS ``Case ID``, 
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS ``Study Code``
"@

$ws.Range("B4").Value = $newQuery

# The author re-selected the edited query cell (and scrolled it into view)
# before saving: selection moves from B3 to B4.
$ws.Range("B4").Select()
$excel.ActiveWindow.ScrollRow = 4
